$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "35.049.76"
$ws.Range("E2").Value = "  +1.18%  "
$ws.Range("D3").Value = "1.844.59"
$ws.Range("E3").Value = "  +1.71%  "
$ws.Range("E4").Value = "  +0.06%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "227.91"
$c.ClearFormats()
$ws.Range("E5").Value = "  +0.76%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "0.614"
$c.ClearFormats()
$ws.Range("E6").Value = "  +2.09%  "
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "1.00"
$c.ClearFormats()
$ws.Range("E7").Value = "  +0.08%  "
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "42.17"
$c.ClearFormats()
$ws.Range("E8").Value = "  +13.66%  "
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "0.0692"
$c.ClearFormats()
$ws.Range("E10").Value = "  +1.06%  "
$ws.Range("E11").Value = "  +3.47%  "
$ws.Range("D12").Value = "2.112.32"
$ws.Range("E12").Value = "  +1.85%  "
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "11.53"
$c.ClearFormats()
$ws.Range("E13").Value = "  +1.38%  "
$ws.Range("D14").Value = "1.836.70"
$ws.Range("E14").Value = "  +0.04%  "
$ws.Range("E15").Value = "  +6.78%  "
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "0.661"
$c.ClearFormats()
$ws.Range("E16").Value = "  +4.05%  "
$ws.Range("D17").Value = "34.919.26"
$ws.Range("E17").Value = "  +0.96%  "
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "69.67"
$c.ClearFormats()
$ws.Range("E18").Value = "  +1.07%  "
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "245.36"
$c.ClearFormats()
$ws.Range("E19").Value = "  +0.60%  "
$ws.Range("D20").Value = "0.0₃0791"
$ws.Range("E20").Value = "  +1.41%  "
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "12.11"
$c.ClearFormats()
$ws.Range("E21").Value = "  +7.42%  "
$ws.Range("E22").Value = "  +15.15%  "
$ws.Range("E23").Value = "  +0.10%  "
$ws.Range("E24").Value = "  -1.43%  "
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "172.21"
$c.ClearFormats()
$ws.Range("E25").Value = "  +0.19%  "
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "7.93"
$c.ClearFormats()
$ws.Range("E26").Value = "  +0.28%  "
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "17.89"
$c.ClearFormats()
$ws.Range("E27").Value = "  +3.47%  "
$ws.Range("E28").Value = "  +0.87%  "
$ws.Range("E29").Value = "  +0.15%  "
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "1.34"
$c.ClearFormats()
$ws.Range("E30").Value = "  +8.29%  "
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "3.96"
$c.ClearFormats()
$ws.Range("E31").Value = "  +3.04%  "
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "4.04"
$c.ClearFormats()
$ws.Range("E32").Value = "  +2.45%  "
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "0.0535"
$c.ClearFormats()
$ws.Range("E33").Value = "  +3.29%  "
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "1.93"
$c.ClearFormats()
$ws.Range("E34").Value = "  +5.66%  "
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "90.78"
$c.ClearFormats()
$ws.Range("E35").Value = "  +11.33%  "
$ws.Range("E36").Value = "  +2.20%  "
$ws.Range("E37").Value = "  +3.19%  "
$ws.Range("D38").Value = "1.345.59"
$ws.Range("E38").Value = "  -1.56%  "
$ws.Range("E39").Value = "  +9.36%  "
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "2.43"
$c.ClearFormats()
$ws.Range("E40").Value = "  +2.23%  "
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "0.0193"
$c.ClearFormats()
$ws.Range("E41").Value = "  +3.03%  "
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "14.96"
$c.ClearFormats()
$ws.Range("E42").Value = "  +8.75%  "
$ws.Range("E43").Value = "  +6.63%  "
$ws.Range("E44").Value = "  +1.74%  "
$ws.Range("E45").Value = "  +0.80%  "
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "0.0518"
$c.ClearFormats()
$ws.Range("E46").Value = "  +2.86%  "
$ws.Range("E47").Value = "  +3.57%  "
$ws.Range("D48").Value = "2.011.48"
$ws.Range("E48").Value = "  +1.91%  "
$ws.Range("E49").Value = "  +0.09%  "
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "103.45"
$c.ClearFormats()
$ws.Range("B51").Value = "BitcoinSV"
$ws.Range("C51").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "49.13"
$c.ClearFormats()
$ws.Range("E51").Value = "  +1.49%  "
